$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Week 34 (AK) header label ---
# The header row (row 1) stores the week numbers as text ("1".."33"),
# all sharing style index 1 (bold, centered) with the other header cells.
# Force the new header cell to be stored as TEXT (not a number) while reusing
# the same style as its neighbour (AJ1).
$ws.Range("AK1").NumberFormat = "@"
$ws.Range("AK1").Value = "34"
$ws.Range("AJ1").Copy() | Out-Null
$ws.Range("AK1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Week 34 (AK) data values, plus two late-arriving week 33 (AJ) corrections ---
$ws.Range("AK2").Value = 0
$ws.Range("AK3").Value = 0
$ws.Range("AK5").Value = 0
$ws.Range("AK6").Value = 42
$ws.Range("AK7").Value = 2
$ws.Range("AK8").Value = 29
$ws.Range("AK10").Value = 0
$ws.Range("AK11").Value = 0
$ws.Range("AK12").Value = 0
$ws.Range("AK13").Value = 0
$ws.Range("AK14").Value = 0
$ws.Range("AK17").Value = 0
$ws.Range("AK18").Value = 0
$ws.Range("AK19").Value = 0
$ws.Range("AK22").Value = 0
$ws.Range("AK23").Value = 0
$ws.Range("AK24").Value = 0
$ws.Range("AK25").Value = 1
$ws.Range("AK27").Value = 0
$ws.Range("AJ28").Value = 6   # week 33 correction
$ws.Range("AK28").Value = 7
$ws.Range("AK29").Value = 2
$ws.Range("AJ30").Value = 16   # week 33 correction
$ws.Range("AK30").Value = 17
$ws.Range("AK31").Value = 0
$ws.Range("AK34").Value = 0
$ws.Range("AK35").Value = 17
$ws.Range("AK36").Value = 1
$ws.Range("AK37").Value = 0
$ws.Range("AK38").Value = 0
$ws.Range("AK40").Value = 0
$ws.Range("AK41").Value = 0
$ws.Range("AK42").Value = 0
$ws.Range("AK43").Value = 0
$ws.Range("AK44").Value = 0
$ws.Range("AK45").Value = 0
$ws.Range("AK46").Value = 0
$ws.Range("AK47").Value = 0
$ws.Range("AK48").Value = 0
$ws.Range("AK49").Value = 0
$ws.Range("AK50").Value = 0
$ws.Range("AK51").Value = 0
$ws.Range("AK52").Value = 0
$ws.Range("AK53").Value = 0
$ws.Range("AK54").Value = 0
$ws.Range("AK55").Value = 0
$ws.Range("AK56").Value = 0
$ws.Range("AK57").Value = 0
$ws.Range("AK58").Value = 0
